$d = $word.ActiveDocument

# 1. Update wording: "...at least one numeric digit" -> "...at least one digit"
$d.Content.Find.Execute(
    "Password must contain at least one numeric digit",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Password must contain at least one digit", 2) | Out-Null

# 2. The paragraph that used to hold only "    " (whitespace, right after the
#    "...one digit" line) loses its run entirely, becoming a bare empty
#    paragraph. Locate the landmark paragraph by content, then address the
#    following paragraph by numeric Index (".Next" returns an object whose
#    Range bounds don't resolve in this host, so avoid it). Delete only the
#    text characters (not the trailing paragraph mark) so the paragraph
#    itself survives as an empty <w:p/> with no run.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Password must contain at least one digit*") {
        $idx = $p.Index
        $whitespacePara = $d.Paragraphs.Item($idx + 1)
        $textOnly = $d.Range($whitespacePara.Range.Start, $whitespacePara.Range.End - 1)
        $textOnly.Delete() | Out-Null
        break
    }
}

# 3. Right after "return True, ..." there used to be two consecutive blank
#    paragraphs; now there is only one. Delete one of them.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*return True, *Password is valid*") {
        $idx = $p.Index
        $blank1 = $d.Paragraphs.Item($idx + 1)
        $blank1.Range.Delete() | Out-Null
        break
    }
}

# 4. A new "  print(valid, msg)" paragraph is inserted between the two blank
#    paragraphs that precede "if __name__ == "__main__":" (i.e. right after
#    the first of those two blanks, which directly follows the final
#    existing "  print(valid, msg)" line).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*if __name__*__main__*") {
        $idx = $p.Index
        $firstBlank = $d.Paragraphs.Item($idx - 2)
        $firstBlank.Range.InsertParagraphAfter() | Out-Null
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*if __name__*__main__*") {
        $idx = $p.Index
        $newPara = $d.Paragraphs.Item($idx - 2)
        $newPara.Range.Text = "  print(valid, msg)"
        break
    }
}
